# ISP20220288 - 調整
# 1. UI上的 Released Date 要改成 Report Date
# 2. Excel 加上欄位 (Received Date / Report Date)
# 3. 日期的框框寬度加長 (Test Date column shifts from J to L; new J/K get the
#    same width/format as the original date column)
#
# Sheet1 currently has a single header row:
#   A Type | B Report No | C SP# | D Brand | E Style# | F Season |
#   G Article | H Artwork | I Result | J Test Date
#
# Target header row:
#   A Type | B Report No | C SP# | D Brand | E Style# | F Season |
#   G Article | H Artwork | I Result | J Received Date | K Report Date |
#   L Test Date

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original "Test Date" header cell (column J) so its value and
# formatting (date number format, fill, border, font) can be reproduced on
# the new columns without disturbing the sheet's column grid (inserting
# whole columns shifts the trailing default column range past the sheet's
# real boundary in this host, so we move data manually instead of using
# Insert).
$testDateValue = $ws.Range("J1").Value()

# 1. Move "Test Date" from J1 to L1, carrying over its style (date format,
#    shaded fill, border).
$ws.Range("J1").Copy()
$ws.Range("L1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("L1").Value = $testDateValue

# 2. Stamp the same style onto J1:K1, then set their new labels.
$ws.Range("L1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("J1").Value = "Received Date"
$ws.Range("K1").Value = "Report Date"

# 3. Widen the new K/L columns to match the original (wider) date column J.
$dateColWidth = $ws.Columns.Item(10).ColumnWidth()
$ws.Columns.Item(11).ColumnWidth = $dateColWidth
$ws.Columns.Item(12).ColumnWidth = $dateColWidth
